# Apply KPI recomputation changes described in the commit.
$wb = $excel.ActiveWorkbook

# --- Sheet "Productdata": update SetupCosts (col E) and StartingInventories (col C) ---
$wsProd = $wb.Worksheets.Item("Productdata")

$prodE = @{
    2  = 0.61875
    3  = 0.9231750000000001
    4  = 0.9364500000000001
    5  = 0.9556650000000001
    6  = 0.54
    7  = 1.58625
    8  = 0.8894250000000001
    9  = 1.7604
    10 = 0.2475
    11 = 0.54
    12 = 0.36
    13 = 0.7424999999999999
    14 = 1.0125
    15 = 0.54
    16 = 0.3381750000000001
    17 = 0.6758999999999999
    18 = 1.58625
}
foreach ($row in $prodE.Keys) {
    $wsProd.Range("E$row").Value = $prodE[$row]
}

$prodCZeroRows = @(8, 9, 10, 11, 12, 13, 15, 16, 17, 18)
foreach ($row in $prodCZeroRows) {
    $wsProd.Range("C$row").Value = 0
}

# --- Sheet "Capacity": update capacity values (col B) ---
$wsCap = $wb.Worksheets.Item("Capacity")

$capB = @{
    2  = 20
    3  = 10
    4  = 20
    5  = 15
    6  = 5
    7  = 15
    8  = 5
    9  = 50
    10 = 20
    11 = 20
    12 = 60
    13 = 15
    14 = 75
    15 = 10
    16 = 10
    17 = 10
    18 = 30
}
foreach ($row in $capB.Keys) {
    $wsCap.Range("B$row").Value = $capB[$row]
}

# --- Sheet "ProcessingTime": update diagonal processing times ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$procDiag = @{
    2  = @{ Col = "B"; Value = 4 }
    5  = @{ Col = "E"; Value = 3 }
    6  = @{ Col = "F"; Value = 1 }
    7  = @{ Col = "G"; Value = 1 }
    8  = @{ Col = "H"; Value = 1 }
    9  = @{ Col = "I"; Value = 5 }
    10 = @{ Col = "J"; Value = 4 }
    11 = @{ Col = "K"; Value = 1 }
    12 = @{ Col = "L"; Value = 3 }
    13 = @{ Col = "M"; Value = 1 }
    14 = @{ Col = "N"; Value = 5 }
    15 = @{ Col = "O"; Value = 2 }
    16 = @{ Col = "P"; Value = 2 }
    17 = @{ Col = "Q"; Value = 1 }
    18 = @{ Col = "R"; Value = 2 }
}
foreach ($row in $procDiag.Keys) {
    $entry = $procDiag[$row]
    $wsProc.Range("$($entry.Col)$row").Value = $entry.Value
}

$wb.Save()
